# Update cryptos price/volume figures (scraped refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.807.38"
$ws.Range("E2").Value = "  -0.45%  "
$ws.Range("D3").Value = "2.312.99"
$ws.Range("E3").Value = "  +3.18%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'97.11"
$ws.Range("E5").Value = "  +3.53%  "
$ws.Range("D6").Value = "'272.49"
$ws.Range("E6").Value = "  +0.63%  "
$ws.Range("D7").Value = "'0.629"
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "'0.626"
$ws.Range("E9").Value = "  -1.26%  "
$ws.Range("D10").Value = "'45.36"
$ws.Range("E10").Value = "  -1.34%  "
$ws.Range("D11").Value = "'0.0955"
$ws.Range("E11").Value = "  -1.73%  "
$ws.Range("D12").Value = "'7.99"
$ws.Range("E12").Value = "  -3.79%  "
$ws.Range("D13").Value = "'0.106"
$ws.Range("E13").Value = "  +0.71%  "
$ws.Range("D14").Value = "2.649.17"
$ws.Range("E14").Value = "  +2.61%  "
$ws.Range("D15").Value = "'15.50"
$ws.Range("E15").Value = "  +1.82%  "
$ws.Range("E16").Value = "  +7.10%  "
$ws.Range("D17").Value = "2.316.48"
$ws.Range("E17").Value = "  +2.75%  "
$ws.Range("D18").Value = "43.753.46"
$ws.Range("E18").Value = "  -0.49%  "
$ws.Range("E19").Value = "  +4.16%  "
$ws.Range("E20").Value = "  +4.39%  "
$ws.Range("D21").Value = "'73.46"
$ws.Range("E21").Value = "  +3.69%  "
$ws.Range("D22").Value = "'239.51"
$ws.Range("E22").Value = "  +1.92%  "
$ws.Range("E23").Value = "  -2.22%  "
$ws.Range("D24").Value = "'9.47"
$ws.Range("E24").Value = "  +3.46%  "
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("E26").Value = "  +1.29%  "
$ws.Range("D27").Value = "'11.37"
$ws.Range("E27").Value = "  -0.52%  "
$ws.Range("E28").Value = "  -2.10%  "
$ws.Range("D29").Value = "'2.29"
$ws.Range("E29").Value = "  +1.31%  "
$ws.Range("D30").Value = "'38.25"
$ws.Range("E30").Value = "  -6.20%  "
$ws.Range("D31").Value = "'22.44"
$ws.Range("E31").Value = "  +7.03%  "
$ws.Range("D32").Value = "'175.12"
$ws.Range("E32").Value = "  +1.45%  "
$ws.Range("D33").Value = "'0.0914"
$ws.Range("E33").Value = "  -0.15%  "
$ws.Range("D34").Value = "'5.49"
$ws.Range("E34").Value = "  +0.15%  "
$ws.Range("E35").Value = "  +2.13%  "
$ws.Range("E36").Value = "  +2.97%  "
$ws.Range("E37").Value = "  -4.22%  "
$ws.Range("D38").Value = "'4.46"
$ws.Range("E38").Value = "  +2.91%  "
$ws.Range("E39").Value = "  -6.26%  "
$ws.Range("D40").Value = "'0.244"
$ws.Range("E40").Value = "  +8.09%  "
$ws.Range("E41").Value = "  +10.45%  "
$ws.Range("D42").Value = "'1.43"
$ws.Range("E42").Value = "  +23.44%  "
$ws.Range("D43").Value = "'12.36"
$ws.Range("E43").Value = "  -4.58%  "
$ws.Range("D44").Value = "'62.89"
$ws.Range("E44").Value = "  -1.30%  "
$ws.Range("D45").Value = "'9.19"
$ws.Range("E45").Value = "  +9.63%  "
$ws.Range("E46").Value = "  -1.28%  "
$ws.Range("E47").Value = "  +4.01%  "
$ws.Range("D48").Value = "'100.53"
$ws.Range("E48").Value = "  -0.30%  "
$ws.Range("E49").Value = "  +0.27%  "
$ws.Range("E50").Value = "  +14.99%  "
$ws.Range("D51").Value = "2.536.45"
$ws.Range("E51").Value = "  +2.81%  "
